# tambahan detail info penyelesaian proses
# Adds a new "Sheet2" after the existing "Sheet1" with summary statistics
# about the scheduling simulation results (total processes, average /
# total waiting time, average / total turn-around time).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1 so the tab order becomes
# Sheet1, Sheet2 (matching the target workbook).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Total Process"
$ws2.Range("B1").Value = 100

$ws2.Range("A2").Value = "AWT (Average Waiting Time)"
$ws2.Range("B2").Value = 476.13

$ws2.Range("A3").Value = "Total Waiting Time"
$ws2.Range("B3").Value = 47613

$ws2.Range("A4").Value = "ATAT (Average Turn Around Time)"
$ws2.Range("B4").Value = 488.71

$ws2.Range("A5").Value = "Total Turn Around Time"
$ws2.Range("B5").Value = 48871

# Restore Sheet1 as the active tab, since adding a sheet makes the new
# sheet active by default.
$ws1.Activate()
